# tickets_pideu.xlsx -- fixing beta_test.py / v.10.9
#
# Append the 2024-05-16 12:27-12:43 incident batch (rows 48-55) to Sheet1.
# (The pre-existing blank placeholder row 2 carries no data, so it simply
# drops out of the sheet once the workbook round-trips through save -
# nothing further needs to be done for that part of the change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number, then values for columns A..I ($null => leave the cell untouched)
$newRows = @(
    @(48, @('2024-05-16', '12:27:56', 'Fallo tolva', '-', '-', '-', '-', '12:28:01', '0:00:05')),
    @(49, @('2024-05-16', '12:28:00', 'Secuencia atornillador', '-', '-', '-', '-', '12:32:07', '0:04:07')),
    @(50, @('2024-05-16', '12:28:17', 'Palet atascado en la curva', '-', '-', '-', '-', '12:32:28', '0:04:11')),
    @(51, @('2024-05-16', '12:32:45', 'Fallo atornillador', '-', '-', '-', '-', '12:32:52', '0:00:07')),
    @(52, @('2024-05-16', '12:42:26', 'Palet atascado en la curva', '-', '-', '-', '-', $null, $null)),
    @(53, @('2024-05-16', '12:42:29', 'Fallo tolva', '-', '-', '-', '-', $null, $null)),
    @(54, @('2024-05-16', '12:43:36', 'No atornilla tapa', '-', '-', '-', '-', $null, $null)),
    @(55, @('2024-05-16', '12:43:40', 'No pone tornillo', '-', '-', '-', '-', '', ''))
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $v = $vals[$c]
        if ($v -eq $null) { continue }

        $cell = $ws.Cells.Item($r, $c + 1)
        if ($v -eq '2024-05-16') {
            # leading apostrophe forces text so Excel doesn't reinterpret
            # the ISO date-looking string as a real date serial number
            $cell.Value = "'" + $v
        } elseif ($v -eq '') {
            # same trick to store a literal empty text value in H55/I55
            $cell.Value = "'"
        } else {
            $cell.Value = $v
        }
    }
}
